$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '49.981.91'
$ws.Cells.Item(2, 5).Value = '  +3.78%  '
$ws.Cells.Item(3, 4).Value = '2.650.44'
$ws.Cells.Item(3, 5).Value = '  +6.01%  '
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).Value = '113.87'
$ws.Cells.Item(5, 5).Value = '  +7.55%  '
$ws.Cells.Item(6, 4).Value = '325.60'
$ws.Cells.Item(6, 5).Value = '  +2.36%  '
$ws.Cells.Item(7, 5).Value = '  +2.06%  '
$ws.Cells.Item(9, 4).Value = '0.558'
$ws.Cells.Item(9, 5).Value = '  +3.88%  '
$ws.Cells.Item(10, 4).Value = '41.03'
$ws.Cells.Item(10, 5).Value = '  +5.63%  '
$ws.Cells.Item(11, 4).Value = '20.10'
$ws.Cells.Item(11, 5).Value = '  -0.62%  '
$ws.Cells.Item(12, 4).Value = '0.0824'
$ws.Cells.Item(12, 5).Value = '  +2.64%  '
$ws.Cells.Item(13, 5).Value = '  +0.83%  '
$ws.Cells.Item(14, 4).Value = '7.39'
$ws.Cells.Item(14, 5).Value = '  +4.08%  '
$ws.Cells.Item(15, 4).Value = '3.063.68'
$ws.Cells.Item(15, 5).Value = '  +5.89%  '
$ws.Cells.Item(16, 4).Value = '2.639.33'
$ws.Cells.Item(16, 5).Value = '  +5.44%  '
$ws.Cells.Item(17, 4).Value = '0.876'
$ws.Cells.Item(17, 5).Value = '  +5.56%  '
$ws.Cells.Item(18, 4).Value = '49.930.60'
$ws.Cells.Item(18, 5).Value = '  +3.93%  '
$ws.Cells.Item(19, 4).Value = '13.24'
$ws.Cells.Item(19, 5).Value = '  +2.91%  '
$ws.Cells.Item(20, 5).Value = '  +2.48%  '
$ws.Cells.Item(21, 4).Value = '2.93'
$ws.Cells.Item(21, 5).Value = '  -2.44%  '
$ws.Cells.Item(23, 4).Value = '72.44'
$ws.Cells.Item(23, 5).Value = '  +1.81%  '
$ws.Cells.Item(24, 4).Value = '275.34'
$ws.Cells.Item(24, 5).Value = '  +2.87%  '
$ws.Cells.Item(25, 4).Value = '2.59'
$ws.Cells.Item(25, 5).Value = '  +2.77%  '
$ws.Cells.Item(26, 4).Value = '26.87'
$ws.Cells.Item(26, 5).Value = '  +4.19%  '
$ws.Cells.Item(28, 4).Value = '10.02'
$ws.Cells.Item(28, 5).Value = '  +2.90%  '
$ws.Cells.Item(29, 5).Value = '  -1.37%  '
$ws.Cells.Item(30, 4).Value = '36.55'
$ws.Cells.Item(30, 5).Value = '  +5.89%  '
$ws.Cells.Item(31, 4).Value = '0.142'
$ws.Cells.Item(31, 5).Value = '  +2.38%  '
$ws.Cells.Item(32, 4).Value = '50.22'
$ws.Cells.Item(32, 5).Value = '  +1.68%  '
$ws.Cells.Item(33, 4).Value = '5.49'
$ws.Cells.Item(33, 5).Value = '  +3.76%  '
$ws.Cells.Item(34, 4).Value = '0.0816'
$ws.Cells.Item(34, 5).Value = '  +5.57%  '
$ws.Cells.Item(35, 4).Value = '19.46'
$ws.Cells.Item(35, 5).Value = '  +1.60%  '
$ws.Cells.Item(36, 5).Value = '  -0.11%  '
$ws.Cells.Item(37, 5).Value = '  +9.91%  '
$ws.Cells.Item(38, 5).Value = '  +6.37%  '
$ws.Cells.Item(39, 5).Value = '  +8.08%  '
$ws.Cells.Item(40, 4).Value = '124.25'
$ws.Cells.Item(40, 5).Value = '  +0.98%  '
$ws.Cells.Item(41, 5).Value = '  +2.10%  '
$ws.Cells.Item(42, 5).Value = '  +0.19%  '
$ws.Cells.Item(43, 4).Value = '22.15'
$ws.Cells.Item(43, 5).Value = '  -2.03%  '
$ws.Cells.Item(44, 4).Value = '0.0320'
$ws.Cells.Item(44, 5).Value = '  +5.80%  '
$ws.Cells.Item(45, 4).Value = '2.083.89'
$ws.Cells.Item(45, 5).Value = '  +4.15%  '
$ws.Cells.Item(46, 5).Value = '  +6.08%  '
$ws.Cells.Item(47, 5).Value = '  +13.63%  '
$ws.Cells.Item(48, 5).Value = '  +4.46%  '
$ws.Cells.Item(49, 4).Value = '9.15'
$ws.Cells.Item(49, 5).Value = '  +2.29%  '
$ws.Cells.Item(50, 2).Value = 'MultiversX'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(50, 4).Value = '60.60'
$ws.Cells.Item(50, 5).Value = '  +7.31%  '
$ws.Cells.Item(51, 2).Value = 'THORChain'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(51, 4).Value = '5.34'
$ws.Cells.Item(51, 5).Value = '  +2.63%  '
